$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they are not
# auto-converted to numbers (matches the source data which stores these
# as literal text, e.g. "149.30" must stay "149.30", not become 149.3).
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D17","D18","D19","D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.153.72"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.899.42"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "306.93"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "0.5234"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "0.3807"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "0.07288"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "21.41"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "0.9032"
$ws.Range("D12").Value = "0.08164"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").Value = "95.37"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.344"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.809.67"
$ws.Range("E15").Value = "  -4.85%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "0.000008644"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "14.69"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "27.187.55"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "5.099"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "10.79"
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").Value = "6.451"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "2.325"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").Value = "149.30"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "1.742"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "115.79"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "4.825"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "4.881"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "0.09218"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").Value = "0.05048"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "0.7921"
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").Value = "1.224"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "2.983"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").Value = "3.365"
$ws.Range("D37").Value = "2.651"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "0.5710"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "1.079"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").Value = "9.031"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "6.590"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "116.29"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "0.1513"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "0.4888"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "10.12"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "38.40"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").Value = "63.89"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "0.05956"
